$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new applicant data
$ws.Range("A2").Value = "SANTIAGO"
$ws.Range("B2").Value = "RAMIREZ"
$ws.Range("C2").Value = "VALENCIA"
$ws.Range("D2").Value = "T.I."
$ws.Range("E2").Value = 1001
$ws.Range("F2").Value = "CARRERA 64"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "01-09-2001"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 23
$ws.Range("I2").Value = "CASTILLA"
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 2014567899
$ws.Range("L2").Value = "ramirez12@gmail.com"
$ws.Range("M2").Value = "Tecnología en Gestión de Mercadeo - Sede Robledo"
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = "uploaded_files\1001_CÉDULA.pdf"
$ws.Range("P2").Value = "uploaded_files\1001_CIVICA.pdf"
$ws.Range("Q2").Value = "uploaded_files\1001_SERVICIOPUBLICOS.pdf"
$ws.Range("R2").Value = "uploaded_files\1001_ANEXO1.pdf"
$ws.Range("S2").Value = "uploaded_files\1001_ANEXO2.xlsx"

# Remove row 3 entirely (it was deleted in the new version)
$ws.Rows("3").Delete()
